$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "243.80"
$ws.Range("E2").Value = "-0.10%"
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "30.00"
$ws.Range("E3").Value = "13.62%"
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.154"
$ws.Range("E4").Value = "-0.08%"
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.05675"
$ws.Range("E5").Value = "1.25%"
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "6.541"
$ws.Range("E6").Value = "1.13%"
$rng.Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "0.8486"
$ws.Range("E7").Value = "3.66%"
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "0.8615"
$ws.Range("E8").Value = "4.05%"
$rng.Style = "Normal"

$rng = $ws.Range("B9:E9")
$rng.NumberFormat = "@"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "0.01008"
$ws.Range("E9").Value = "0.43%"
$rng.Style = "Normal"

$rng = $ws.Range("B10:E10")
$rng.NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1346"
$ws.Range("E10").Value = "1.05%"
$rng.Style = "Normal"

$rng = $ws.Range("B11:E11")
$rng.NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.06907"
$ws.Range("E11").Value = "-0.48%"
$rng.Style = "Normal"

$rng = $ws.Range("B12:E12")
$rng.NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.02890"
$ws.Range("E12").Value = "-0.19%"
$rng.Style = "Normal"

$rng = $ws.Range("B13:E13")
$rng.NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09372"
$ws.Range("E13").Value = "-0.14%"
$rng.Style = "Normal"

$rng = $ws.Range("B14:E14")
$rng.NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001526"
$ws.Range("E14").Value = "0.28%"
$rng.Style = "Normal"

$rng = $ws.Range("B15:E15")
$rng.NumberFormat = "@"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "0.04174"
$ws.Range("E15").Value = "-9.39%"
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "0.006131"
$ws.Range("E16").Value = "-1.86%"
$rng.Style = "Normal"

$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$ws.Range("E17").Value = "-4.03%"
$rng.Style = "Normal"

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = "3.033"
$ws.Range("E18").Value = "-0.04%"
$rng.Style = "Normal"

$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = "2.243"
$ws.Range("E19").Value = "2.76%"
$rng.Style = "Normal"

$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$ws.Range("E20").Value = "1.19%"
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "0.03351"
$ws.Range("E21").Value = "9.36%"
$rng.Style = "Normal"

$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$ws.Range("E22").Value = "0.31%"
$rng.Style = "Normal"

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "3.630"
$ws.Range("E23").Value = "-3.01%"
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.1374"
$ws.Range("E24").Value = "2.39%"
$rng.Style = "Normal"

$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$ws.Range("E25").Value = "-1.00%"
$rng.Style = "Normal"

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.004439"
$ws.Range("E26").Value = "-1.07%"
$rng.Style = "Normal"

$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$ws.Range("E28").Value = "-0.55%"
$rng.Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.03767"
$ws.Range("E40").Value = "3.45%"
$rng.Style = "Normal"

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.005778"
$ws.Range("E41").Value = "-6.31%"
$rng.Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.1057"
$ws.Range("E42").Value = "0.58%"
$rng.Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.002289"
$ws.Range("E43").Value = "-4.60%"
$rng.Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.009273"
$ws.Range("E44").Value = "2.60%"
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.00005102"
$ws.Range("E45").Value = "-4.54%"
$rng.Style = "Normal"

$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$ws.Range("E46").Value = "0.01%"
$rng.Style = "Normal"

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = "0.08997"
$ws.Range("E47").Value = "-37.49%"
$rng.Style = "Normal"

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.002769"
$ws.Range("E48").Value = "-5.65%"
$rng.Style = "Normal"

$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$ws.Range("E49").Value = "0.01%"
$rng.Style = "Normal"

$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
$rng.Style = "Normal"
